# updated temp data and size data
#
# Appends 18 new temperature-log rows (2025-07-02, T0 and T1 timepoints) to
# the bottom of Sheet1, then updates the sheet's scroll/selection state to
# where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prime the shared-string table -----------------------------------------
# The plate-name strings in column B must be interned (added to
# xl/sharedStrings.xml) in the exact order the author first typed/pasted
# them, which is NOT strictly top-to-bottom row order (plate42 was entered
# before plate41). Touch column B once, in that precise sequence, before
# writing the rows out in full below (re-setting the same value afterwards
# is a no-op for the shared-string table since the string is already interned).
$primeOrder = @(
    @(57, "plate33"),
    @(58, "plate35"),
    @(59, "plate37"),
    @(60, "plate40"),
    @(62, "plate42"),
    @(61, "plate41"),
    @(63, "plate43"),
    @(64, "plate44"),
    @(65, "plate45"),
    @(66, "plate47"),
    @(69, "plate29"),
    @(73, "plate39")
)
foreach ($entry in $primeOrder) {
    $ws.Cells.Item($entry[0], 2).Value = $entry[1]
}

# --- Write the new rows in full (date, plate, timepoint, temperature) ------
$rows = @(
    @(57, "plate33", "T0", 18.4),
    @(58, "plate35", "T0", 16.8),
    @(59, "plate37", "T0", 11.1),
    @(60, "plate40", "T0", 12.2),
    @(61, "plate41", "T0", 9),
    @(62, "plate42", "T0", 7.7),
    @(63, "plate43", "T0", 10.2),
    @(64, "plate44", "T0", 9.9),
    @(65, "plate45", "T0", 9.5),
    @(66, "plate47", "T0", 11.2),
    @(67, "plate25", "T1", 35.4),
    @(68, "plate27", "T1", 32.4),
    @(69, "plate29", "T1", 34.8),
    @(70, "plate31", "T1", 33.8),
    @(71, "plate35", "T1", 34.7),
    @(72, "plate37", "T1", 33.9),
    @(73, "plate39", "T1", 33.7),
    @(74, "plate41", "T1", 33.5)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 20250702
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# --- Match the author's final view state ------------------------------------
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E79").Select()
